# Apply the 2024-09-19 cryptos-list data refresh.
# Each changed cell is set via Range.Value; numeric-looking price strings
# (column D) are written with a leading apostrophe so Excel keeps them as
# text (matching the original "inlineStr" cell type) instead of coercing
# them to Number, and the quote-prefix style that introduces is reset back
# to "Normal" immediately after so no new cell style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.013.54"
$ws.Range("E2").Value = "  +2.60%  "
# Row 3
$ws.Range("D3").Value = "2.417.31"
$ws.Range("E3").Value = "  +3.88%  "
# Row 5
$ws.Range("D5").Value = "'558.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.47%  "
# Row 6
$ws.Range("D6").Value = "'138.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.66%  "
# Row 7
$ws.Range("E7").Value = "  +0.04%  "
# Row 8
$ws.Range("E8").Value = "  +0.70%  "
# Row 9
$ws.Range("D9").Value = "2.415.81"
$ws.Range("E9").Value = "  +3.93%  "
# Row 10
$ws.Range("E10").Value = "  +3.04%  "
# Row 11
$ws.Range("D11").Value = "'5.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.70%  "
# Row 13
$ws.Range("D13").Value = "'0.346"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.46%  "
# Row 14
$ws.Range("D14").Value = "'25.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.08%  "
# Row 15
$ws.Range("D15").Value = "2.848.41"
$ws.Range("E15").Value = "  +3.94%  "
# Row 16
$ws.Range("D16").Value = "61.982.35"
$ws.Range("E16").Value = "  +2.58%  "
# Row 17
$ws.Range("E17").Value = "  +4.80%  "
# Row 18
$ws.Range("D18").Value = "2.416.85"
$ws.Range("E18").Value = "  +4.55%  "
# Row 19
$ws.Range("D19").Value = "'11.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.66%  "
# Row 20
$ws.Range("D20").Value = "'343.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.18%  "
# Row 21
$ws.Range("D21").Value = "'4.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.05%  "
# Row 22
$ws.Range("D22").Value = "'6.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.69%  "
# Row 23
$ws.Range("E23").Value = "  +0.13%  "
# Row 24
$ws.Range("D24").Value = "'64.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "
# Row 25
$ws.Range("E25").Value = "  -0.70%  "
# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
# Row 27
$ws.Range("D27").Value = "'8.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.32%  "
# Row 28
$ws.Range("D28").Value = "'1.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.17%  "
# Row 29
$ws.Range("D29").Value = "'1.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.78%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  +7.14%  "
# Row 31
$ws.Range("D31").Value = "'1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.79%  "
# Row 32
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.35%  "
# Row 33
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'171.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.40%  "
# Row 34
$ws.Range("D34").Value = "'1.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.58%  "
# Row 35
$ws.Range("D35").Value = "'0.394"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.53%  "
# Row 36
$ws.Range("D36").Value = "'374.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.10%  "
# Row 37
$ws.Range("D37").Value = "'18.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.72%  "
# Row 38
$ws.Range("E38").Value = "  +9.99%  "
# Row 40
$ws.Range("E40").Value = "  -0.06%  "
# Row 41
$ws.Range("D41").Value = "'1.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.21%  "
# Row 42
$ws.Range("D42").Value = "'39.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.88%  "
# Row 43
$ws.Range("D43").Value = "'145.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.45%  "
# Row 44
$ws.Range("D44").Value = "'3.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.56%  "
# Row 45
$ws.Range("D45").Value = "'20.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.84%  "
# Row 46
$ws.Range("D46").Value = "'0.0955"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "
# Row 47
$ws.Range("D47").Value = "'0.586"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.04%  "
# Row 48
$ws.Range("D48").Value = "'0.0516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.32%  "
# Row 49
$ws.Range("D49").Value = "'17.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.16%  "
# Row 50
$ws.Range("E50").Value = "  +3.12%  "
# Row 51
$ws.Range("D51").Value = "0.0₆0220"
$ws.Range("E51").Value = "  +3.10%  "
